$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The r / p-value columns hold numeric-looking values that must remain TEXT
# (the sheet stores them as strings, not numbers), so they are entered with a
# leading apostrophe, exactly like a user forcing text entry in Excel.

# Row 2: Cognitive Failures -> Trust
$ws.Range("B2").Value = "'0.067"
$ws.Range("C2").Value = "'0.5350"

# Row 3: Cognitive Load -> Trust
$ws.Range("B3").Value = "'-0.012"
$ws.Range("C3").Value = "'0.9123"

# Row 4: Performance Expectancy -> Trust
$ws.Range("B4").Value = "'0.622"

# Row 5: Effort Expectancy -> Trust
$ws.Range("B5").Value = "'0.626"

# Row 6: Facilitating Conditions -> Trust
$ws.Range("B6").Value = "'0.257"
$ws.Range("C6").Value = "'0.0151"

# Row 7: Cognitive Failures <-> Cognitive Load
$ws.Range("B7").Value = "'0.596"

# Row 8: Cognitive Load -> Effort Expectancy
$ws.Range("B8").Value = "'-0.276"
$ws.Range("C8").Value = "'0.0088"

# Row 9: AI Use Frequency -> Trust
$ws.Range("B9").Value = "'0.354"
$ws.Range("C9").Value = "'0.0007"
$ws.Range("D9").Value = "***"

# Row 10: Education -> Trust
$ws.Range("B10").Value = "'0.016"
$ws.Range("C10").Value = "'0.8809"
